$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "jhasbdjh"
$ws.Range("B5").Value = "asbdjhsa"
$ws.Range("C5").Value = 268361723
$ws.Range("D5").Value = "shadbjsa223"
$ws.Range("E2").Copy()
$ws.Range("E5").PasteSpecial()
$ws.Range("F5").Value = "asbdsad@elpdjcn.dsbh.com"
$ws.Range("G5").Value = "No tiene"
